$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Program1.txt")
$ws2 = $wb.Worksheets.Item("Program 1 conversion")

$ws1.Range("A60").Value = "SOB 1,2,0,1(aifls)"
$ws2.Range("A61").Value = "SOB 1,2,0,1(aifls)"
